$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# Update Price (column D) values
Set-TextValue $ws.Range("D2") '278.76'
Set-TextValue $ws.Range("D3") '27.42'
Set-TextValue $ws.Range("D4") '4.843'
Set-TextValue $ws.Range("D5") '0.06335'
Set-TextValue $ws.Range("D6") '6.950'
Set-TextValue $ws.Range("D7") '3.405'
Set-TextValue $ws.Range("D8") '0.8750'
Set-TextValue $ws.Range("D9") '0.9559'
Set-TextValue $ws.Range("D10") '0.1470'
Set-TextValue $ws.Range("D11") '0.05124'
Set-TextValue $ws.Range("D12") '0.07331'
Set-TextValue $ws.Range("D13") '0.03104'
Set-TextValue $ws.Range("D14") '0.09068'
Set-TextValue $ws.Range("D15") '0.001563'
Set-TextValue $ws.Range("D16") '0.0006288'
Set-TextValue $ws.Range("D17") '0.006081'
Set-TextValue $ws.Range("D18") '3.447'
Set-TextValue $ws.Range("D20") '0.3144'
Set-TextValue $ws.Range("D21") '0.1312'
Set-TextValue $ws.Range("D22") '3.867'
Set-TextValue $ws.Range("D23") '0.04332'
Set-TextValue $ws.Range("D24") '0.001178'
Set-TextValue $ws.Range("D25") '0.004296'
Set-TextValue $ws.Range("D27") '0.0001691'
Set-TextValue $ws.Range("D40") '0.04086'
Set-TextValue $ws.Range("D41") '0.006684'
Set-TextValue $ws.Range("D42") '0.1163'
Set-TextValue $ws.Range("D44") '0.01309'
Set-TextValue $ws.Range("D45") '0.00005218'
Set-TextValue $ws.Range("D47") '2.379'
Set-TextValue $ws.Range("D48") '0.02251'

# Update Volume(1h) (column E) values
Set-TextValue $ws.Range("E2") '6.60%'
Set-TextValue $ws.Range("E3") '2.52%'
Set-TextValue $ws.Range("E4") '3.04%'
Set-TextValue $ws.Range("E5") '2.33%'
Set-TextValue $ws.Range("E6") '3.19%'
Set-TextValue $ws.Range("E7") '7.34%'
Set-TextValue $ws.Range("E8") '2.88%'
Set-TextValue $ws.Range("E9") '4.70%'
Set-TextValue $ws.Range("E10") '4.69%'
Set-TextValue $ws.Range("E11") '0.48%'
Set-TextValue $ws.Range("E12") '3.26%'
Set-TextValue $ws.Range("E13") '-0.11%'
Set-TextValue $ws.Range("E14") '0.28%'
Set-TextValue $ws.Range("E15") '2.18%'
Set-TextValue $ws.Range("E16") '1.87%'
Set-TextValue $ws.Range("E17") '1.93%'
Set-TextValue $ws.Range("E18") '0.02%'
Set-TextValue $ws.Range("E19") '4.76%'
Set-TextValue $ws.Range("E20") '2.34%'
Set-TextValue $ws.Range("E21") '0.18%'
Set-TextValue $ws.Range("E22") '-6.05%'
Set-TextValue $ws.Range("E23") '1.96%'
Set-TextValue $ws.Range("E24") '-0.29%'
Set-TextValue $ws.Range("E25") '5.98%'
Set-TextValue $ws.Range("E27") '3.09%'
Set-TextValue $ws.Range("E40") '3.06%'
Set-TextValue $ws.Range("E41") '61.88%'
Set-TextValue $ws.Range("E42") '4.48%'
Set-TextValue $ws.Range("E44") '-1.42%'
Set-TextValue $ws.Range("E45") '1.04%'
Set-TextValue $ws.Range("E47") '853.79%'
Set-TextValue $ws.Range("E48") '-33.85%'

# Update Hora (column G) from 10 to 11 for all data rows (2-51)
foreach ($row in 2..51) {
    Set-TextValue $ws.Range("G$row") '11'
}
